# Apply data updates to the "Inscricoes" sheet.
# These changes correspond to incremental increases in enrollment counts
# (Inscritos / Pagos / Inscrições homologadas) for a handful of rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

# Row 2: Inscritos (E2) 31 -> 32
$ws.Range("E2").Value = 32

# Row 3: Pagos (F3) 16 -> 17 ; Inscrições homologadas (H3) 19 -> 20
$ws.Range("F3").Value = 17
$ws.Range("H3").Value = 20

# Row 4: Inscritos (E4) 25 -> 26 ; Pagos (F4) 9 -> 10 ; Inscrições homologadas (H4) 12 -> 13
$ws.Range("E4").Value = 26
$ws.Range("F4").Value = 10
$ws.Range("H4").Value = 13

# Row 12: Inscritos (E12) 35 -> 36
$ws.Range("E12").Value = 36

# Row 16: Pagos (F16) 97 -> 98 ; Inscrições homologadas (H16) 185 -> 186
$ws.Range("F16").Value = 98
$ws.Range("H16").Value = 186

# Row 18: Inscritos (E18) 100 -> 101
$ws.Range("E18").Value = 101
